# Update cryptos list values per latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.720.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.132.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +11.23%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("E6").Value = "  -3.76%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.88"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0743"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.74%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.434.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.843"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.125.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +10.63%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.809.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  -7.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +16.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.46%  "
$ws.Range("E29").Value = "  -6.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +62.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.124"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0959"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +12.62%  "
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +17.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.953"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.86%  "
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("E40").Value = "  -8.80%  "
$ws.Range("E41").Value = "  +8.84%  "
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.18%  "
$ws.Range("E44").Value = "  +12.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.360.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0843"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.324.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.81%  "
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  -2.44%  "
